$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths
$ws.Columns.Item(1).ColumnWidth = 94.21875
$ws.Columns.Item(2).ColumnWidth = 185.77734375

# --- Test: Standardablauf ---
# (rows 2-4 were filled in before the A1 header, matching original authoring order)
$ws.Range("A2").Value = "Der Nutzer klickt auf das DropDown-Menü für die Auswahl von Spieler 1."
$ws.Range("B2").Value = "Es wird im DropDown-Menü eine List aller möglichen Spieler angezeigt."

$ws.Range("A3").Value = 'Der Nutzer wählt die Option "Mensch" im DropDown-Menü aus.'
$ws.Range("B3").Value = "Als Spieler 1 wird ein Mensch festgelegt. Im DropDown-Menü wird angezeigt, dass die Option ausgewählt wurde."

$ws.Range("A4").Value = "Der Nutzer klickt auf das DropDown-Menü für die Auswahl von Spieler 2."
$ws.Range("B4").Value = "Es wird im DropDown-Menü eine List aller möglichen Spieler angezeigt."

$ws.Range("A1").Value = "Test: Standardablauf"
$ws.Range("A1").Font.Bold = $true

$ws.Range("A5").Value = 'Der Nutzer wählt die Option "KI 1" im DropDown-Menü aus.'
$ws.Range("B5").Value = 'Als Spieler 1 wird die "KI 1" festgelegt. Im DropDown-Menü wird angezeigt, dass die Option ausgewählt wurde.'

$ws.Range("A6").Value = 'Der Nutzer klickt auf das Feld "Spiel starten".'
$ws.Range("B6").Value = "Es wird in die Spielansicht gewechselt. Das Spielfeld ist leer. Der Graph zeigt das leere Feld und alle möglichen Folgezustände an."

$ws.Range("A7").Value = "Der Nutzer klickt auf ein Feld des Spielfeldes."
$ws.Range("B7").Value = 'Das angeklickte Feld wird mit dem Zeichen "X" versehen. Im Graph wird jetzt der Verlauf plus alle weiteren möglichen Folgezustände angezeigt.'

$ws.Range("A8").Value = 'Der Nutzer klickt auf den "Play-Pfeil".'
$ws.Range("B8").Value = 'Die KI macht einen zufälligen Zug und belegt somit ein zufälliges unbesetztes Feld. Dieses wird mit dem Zeichen "O" versehen. Im Graph wird der erneuerte Verlauf plus alle möglichen Folgezustände angezeigt.'

$ws.Range("A9").Value = 'Der Nutzer klickt abwechselnd auf ein freies Spielfeld und auf den "Play-Pfeil" bis das Spiel vorbei ist.'
$ws.Range("B9").Value = "Das Spielfeld wird laufend aktualisiert. Der Graph zeigt jeweils immer den Verlauf plus alle weiteren möglichen Folgezustände an. Sobald das Spiel vorbei ist, wird das Spielergebnis angezeigt."

$ws.Range("A10").Value = 'Der Nutzer klickt auf den Button "Neustart".'
$ws.Range("B10").Value = "Die Spielerauswahl wird angezeigt."

$ws.Range("A11").Value = 'Der Nutzer behält die Konfiguration bei und klickt auf "Spiel starten"'
$ws.Range("B11").Value = "Es wird in die Spielansicht gewechselt. Das Spielfeld ist leer. Der Graph zeigt das leere Feld und alle möglichen Folgezustände an."

# --- Test: AutoPlay/Pause ---
$ws.Range("A13").Value = "Test: AutoPlay/Pause"
$ws.Range("A13").Font.Bold = $true

$ws.Range("A14").Value = 'Der Nutzer klickt auf das "Play"-Symbol'
$ws.Range("B14").Value = 'Das "Play"-Symbol wird durch ein "Pause"-Symbol ersetzt.'

$ws.Range("A15").Value = 'Der Nutzer klickt auf das "NextMove"-Symbol'
$ws.Range("B15").Value = 'Das "Pause"-Symbol wird durch ein "Play"-Symbol ersetzt. Es wird kein Zug ausgeführt, da der Mensch am Zug ist.'

# Selection / view state
[void]$ws.Range("B16").Select()
$excel.ActiveWindow.Zoom = 100
